# Update the USD Amount (column T) figure on the data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("T2").Value = 516522
